# Rename the inline pictures in the document's two headers and two
# footers so that:
#   - the Pearson logo pictures (in the footers) go from "image2.png" to
#     "image1.png"
#   - the BTEC logo pictures (in the headers) go from "image1.jpg" to
#     "image2.jpg"
# Both the <wp:docPr name="..."/> and the <pic:cNvPr name="..."/> that sit
# inside the drawing need to be updated.
#
# Setting InlineShape.Name only updates <wp:docPr>, and (in this runtime)
# throws for shapes that live in a footer story, so instead we rebuild the
# whole drawing run via Range.InsertXML (Flat-OPC wrapped) and drop the
# stale original picture afterwards.

$d = $word.ActiveDocument

function Update-LogoName($HFRange, $DocPrId, $NewName) {
    # $HFRange : Headers.Item(n) / Footers.Item(n)
    # $DocPrId : original wp:docPr/@id to preserve (COM does not expose it)
    # $NewName : new value for both the wp:docPr and pic:cNvPr "name" attrs

    $shp = $HFRange.Range.InlineShapes.Item(1)
    $descr = $shp.AlternativeText

    # Pull geometry straight from the shape so the re-created drawing
    # matches the original exactly except for the two name attributes.
    $cx = [int]($shp.Width * 914400 / 72)
    $cy = [int]($shp.Height * 914400 / 72)

    $rng = $shp.Range

    $drawing = "<w:r><w:drawing><wp:inline distB=`"0`" distT=`"0`" distL=`"0`" distR=`"0`"><wp:extent cx=`"$cx`" cy=`"$cy`"/><wp:effectExtent b=`"0`" l=`"0`" r=`"0`" t=`"0`"/><wp:docPr descr=`"$descr`" id=`"$DocPrId`" name=`"$NewName`"/><a:graphic><a:graphicData uri=`"http://schemas.openxmlformats.org/drawingml/2006/picture`"><pic:pic><pic:nvPicPr><pic:cNvPr descr=`"$descr`" id=`"0`" name=`"$NewName`"/><pic:cNvPicPr preferRelativeResize=`"0`"/></pic:nvPicPr><pic:blipFill><a:blip r:embed=`"rId1`"/><a:srcRect b=`"0`" l=`"0`" r=`"0`" t=`"0`"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x=`"0`" y=`"0`"/><a:ext cx=`"$cx`" cy=`"$cy`"/></a:xfrm><a:prstGeom prst=`"rect`"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>"

    $xml = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" xmlns:wp=`"http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing`" xmlns:a=`"http://schemas.openxmlformats.org/drawingml/2006/main`" xmlns:pic=`"http://schemas.openxmlformats.org/drawingml/2006/picture`" xmlns:r=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships`"><w:body><w:p>$drawing</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

    $rng.InsertXML($xml)

    # Drop the now-duplicated, stale original picture (it is first in
    # document order; the freshly inserted one comes after it).
    $HFRange.Range.InlineShapes.Item(1).Delete()
}

$sec = $d.Sections.Item(1)

# NOTE: Word's Footers.Item(1)/Headers.Item(1) address the "default"
# (primary) header/footer part, while Item(2) addresses the "first page"
# variant - this was confirmed empirically against the underlying parts,
# and does not follow the physical footer1.xml/footer2.xml file-name
# order. The wp:docPr ids below are the ones already present in each
# part, preserved as-is (only the "name" attribute is being renamed).
Update-LogoName $sec.Footers.Item(1) 4 "image1.png"
Update-LogoName $sec.Footers.Item(2) 2 "image1.png"
Update-LogoName $sec.Headers.Item(1) 3 "image2.jpg"
Update-LogoName $sec.Headers.Item(2) 1 "image2.jpg"
